# Apply the commit's changes to the Transactions sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

# B4: TransCost value changed from 97 to 5.5
$ws.Range("B4").Value = 5.5

# C4/D4 (FromName/ToName) swapped: was Neady/Alan, now Alan/Neady
$ws.Range("C4").Value = "Alan"
$ws.Range("D4").Value = "Neady"

# C5/D5 (FromName/ToName) swapped: was Alan/Neady, now Neady/Alan
$ws.Range("C5").Value = "Neady"
$ws.Range("D5").Value = "Alan"

# Force a full recalculation of all formula cells: the engine's incremental
# dependency tracking does not always pick up changes on cells referenced
# only in the untaken branch of a nested IF, so re-apply every formula to
# force it to be freshly (re)evaluated, then recalc the whole workbook.
foreach ($sheet in $wb.Worksheets) {
    $used = $sheet.UsedRange
    foreach ($cell in $used.Cells) {
        if ($cell.HasFormula) {
            $cell.Formula = $cell.Formula
        }
    }
}
$excel.CalculateFullRebuild()

# Update the active cell selection to B4 (was B5)
$ws.Activate()
$ws.Range("B4").Select()

$wb.Save()
